$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stocks")

# Refresh the return/risk figures for the existing 10 stocks (B: Rendement moyen, C: Risque)
$ws.Cells.Item(2, 2).Value = 0.001054800486479152
$ws.Cells.Item(2, 3).Value = 0.01840606280054741
$ws.Cells.Item(3, 2).Value = 0.0006780267961571966
$ws.Cells.Item(3, 3).Value = 0.02024111117386716
$ws.Cells.Item(4, 2).Value = 0.0003610681506591856
$ws.Cells.Item(4, 3).Value = 0.01465004504620969
$ws.Cells.Item(5, 2).Value = 0.0007868466449651378
$ws.Cells.Item(5, 3).Value = 0.01532730165393188
$ws.Cells.Item(6, 2).Value = 0.001087198475546731
$ws.Cells.Item(6, 3).Value = 0.0181208886644603
$ws.Cells.Item(7, 2).Value = 0.0006999207247814089
$ws.Cells.Item(7, 3).Value = 0.01338750310468702
$ws.Cells.Item(8, 2).Value = 0.0007715763614647289
$ws.Cells.Item(8, 3).Value = 0.0264911297135664
$ws.Cells.Item(9, 2).Value = 0.0007572259421507938
$ws.Cells.Item(9, 3).Value = 0.02257839309019112
$ws.Cells.Item(10, 2).Value = 0.001284666371307356
$ws.Cells.Item(10, 3).Value = 0.01715166976866947
$ws.Cells.Item(11, 2).Value = 0.0005862478355695088
$ws.Cells.Item(11, 3).Value = 0.01705881703527166

# Add 20 more stocks (tickers in column A get the same header/label style as A2:A11)
$ws.Cells.Item(12, 1).Value = "DG"
$ws.Cells.Item(12, 2).Value = 0.0006443425410281341
$ws.Cells.Item(12, 3).Value = 0.01938833452923468
$ws.Cells.Item(13, 1).Value = "CS"
$ws.Cells.Item(13, 2).Value = 0.0007838040427657468
$ws.Cells.Item(13, 3).Value = 0.01768672641101808
$ws.Cells.Item(14, 1).Value = "SAF"
$ws.Cells.Item(14, 2).Value = 0.0007002281797929685
$ws.Cells.Item(14, 3).Value = 0.02536189857201171
$ws.Cells.Item(15, 1).Value = "RI"
$ws.Cells.Item(15, 2).Value = 0.0003004166615604259
$ws.Cells.Item(15, 3).Value = 0.01392035340424598
$ws.Cells.Item(16, 1).Value = "KER"
$ws.Cells.Item(16, 2).Value = 0.0002830575929890021
$ws.Cells.Item(16, 3).Value = 0.02016231370586485
$ws.Cells.Item(17, 1).Value = "STLAM"
$ws.Cells.Item(17, 2).Value = 0.001052548662563955
$ws.Cells.Item(17, 3).Value = 0.02360171653360764
$ws.Cells.Item(18, 1).Value = "BN"
$ws.Cells.Item(18, 2).Value = 0.0001992199499318547
$ws.Cells.Item(18, 3).Value = 0.01331857308442761
$ws.Cells.Item(19, 1).Value = "STMPA"
$ws.Cells.Item(19, 2).Value = 0.001348709788091363
$ws.Cells.Item(19, 3).Value = 0.0255878814756731
$ws.Cells.Item(20, 1).Value = "SGO"
$ws.Cells.Item(20, 2).Value = 0.0009242435870429643
$ws.Cells.Item(20, 3).Value = 0.02101635462672052
$ws.Cells.Item(21, 1).Value = "ENGI"
$ws.Cells.Item(21, 2).Value = 0.0006493756917204952
$ws.Cells.Item(21, 3).Value = 0.01742686089911346
$ws.Cells.Item(22, 1).Value = "CAP"
$ws.Cells.Item(22, 2).Value = 0.0008570599752261347
$ws.Cells.Item(22, 3).Value = 0.01992709773705277
$ws.Cells.Item(23, 1).Value = "DSY"
$ws.Cells.Item(23, 2).Value = 0.0007657714506419176
$ws.Cells.Item(23, 3).Value = 0.01891410346755044
$ws.Cells.Item(24, 1).Value = "LR"
$ws.Cells.Item(24, 2).Value = 0.0006588865500131063
$ws.Cells.Item(24, 3).Value = 0.01613309755763861
$ws.Cells.Item(25, 1).Value = "GLE"
$ws.Cells.Item(25, 2).Value = 0.0005989024886458441
$ws.Cells.Item(25, 3).Value = 0.0281664854906935
$ws.Cells.Item(26, 1).Value = "ML"
$ws.Cells.Item(26, 2).Value = 0.0005547833136926399
$ws.Cells.Item(26, 3).Value = 0.01835299589352088
$ws.Cells.Item(27, 1).Value = "ORA"
$ws.Cells.Item(27, 2).Value = 0.00009254655545388085
$ws.Cells.Item(27, 3).Value = 0.01230540693457924
$ws.Cells.Item(28, 1).Value = "VIE"
$ws.Cells.Item(28, 2).Value = 0.0006787764957448416
$ws.Cells.Item(28, 3).Value = 0.01853844246826758
$ws.Cells.Item(29, 1).Value = "PUB"
$ws.Cells.Item(29, 2).Value = 0.0007104758242159748
$ws.Cells.Item(29, 3).Value = 0.02091540024422309
$ws.Cells.Item(30, 1).Value = "ACA"
$ws.Cells.Item(30, 2).Value = 0.0006662060084863733
$ws.Cells.Item(30, 3).Value = 0.02113920399507774
$ws.Cells.Item(31, 1).Value = "EDEN"
$ws.Cells.Item(31, 2).Value = 0.0005865730821049534
$ws.Cells.Item(31, 3).Value = 0.01738613732998491

# Carry over the ticker-cell style (bold, bordered, centered) to the new rows
$ws.Range("A2").Copy()
$ws.Range("A12:A31").PasteSpecial(-4122)
$excel.CutCopyMode = 0
